$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the values that changed in the "algo" table
$ws.Range("I3").Value = 2
$ws.Range("K3").Value = 8

$ws.Range("G5").Value = 4
$ws.Range("I5").Value = 8

$ws.Range("G7").Value = 4
$ws.Range("I7").Value = 6
$ws.Range("K7").Value = 0

# Move the active cell selection to G6, matching the saved view state
$ws.Activate()
$ws.Range("G6").Select()
